$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Drop the columns that are no longer part of the scraper's export:
#   A  (Norm, Typ)
#   G:M (Ritningsnummer, Position, Beteckning, Kompletterande
#        Information ovrigt, Ref annan, Historiskt Varumarke,
#        Historiskt inkopsreferens)
#   O  (Forpackning)
#   Q  (the always-empty spacer column)
# Deleting right-to-left keeps the earlier letters stable while we work.
# ------------------------------------------------------------------
$ws.Columns("Q").Delete()
$ws.Columns("O").Delete()
$ws.Columns("G:M").Delete()
$ws.Columns("A").Delete()

# The remaining 9 columns (B,C,D,E,F,N,P,R,S in the old numbering) are now
# A:I -> Varumarke, Artikelbenamning, GVM, Artikelnummer, Typbeteckning,
#        Enhet, SSG-notering, E-nummer, RSK-nummer

# Give the "SSG-notering" column (now G) a wider, manually sized column.
$ws.Columns("G").ColumnWidth = 28.6

# Turn the old autofilter off - the table no longer carries one.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# The hidden _FilterDatabase defined name survives the filter removal but
# needs its reference shrunk to the new, narrower table extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$I`$21217"
    }
}

# Reset the view: no more frozen/scrolled-to column J, and the active
# selection now targets the new last column (H, E-nummer) instead of Q.
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H1:H1048576").Select() | Out-Null
